$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The first data row (old row 2, date 2007-06-xx / 39400) is removed; all subsequent
# rows shift up by one. Deleting the row preserves the original stored precision for
# every cell whose value does not actually change.
$ws.Rows.Item(2).Delete()

# The naive forecaster recompute changed a handful of values beyond the plain shift:
# column C row 5 (rounding fix) and the full "y_1_forecast" column E, which holds newly
# computed forecasts (and is now blank for the first four data rows).
$ws.Range("E2").ClearContents()
$ws.Range("E3").ClearContents()
$ws.Range("E4").ClearContents()
$ws.Range("C5").Value = 1.173294700162053
$ws.Range("E5").ClearContents()
$ws.Range("E6").Value = 1.1370912555561
$ws.Range("E7").Value = 0.9849212343369107
$ws.Range("E8").Value = 1.242282657891813
$ws.Range("E9").Value = 1.358051868183585
$ws.Range("E10").Value = 1.377345568933785
$ws.Range("E11").Value = 1.476362359157601
$ws.Range("E12").Value = 1.451677407676555
$ws.Range("E13").Value = 0.6182077276742692
$ws.Range("E14").Value = -1.655311137157178
$ws.Range("E15").Value = 4.631210905746741
$ws.Range("E16").Value = 1.066562775371072
$ws.Range("E17").Value = 0.4316736535407095
$ws.Range("E18").Value = 0.8173856700710358
